# My Courses Distribution Scenario for Available Courses Completed
# - flips the "Forgot Password" suite's Runmode from YES to NO
# - appends a new "MC Distribution" test-suite row
# - moves the active-cell selection down to where the next row would be entered

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# The "Forgot Password" suite is no longer run (YES -> NO).
$ws.Range("C4").Value = "NO"

# New row: My Courses distribution scenario.
$ws.Range("A6").Value = "MC Distribution"
$ws.Range("B6").Value = "My Courses distribution description"
$ws.Range("C6").Value = "Yes"

# Update the sheet's active selection to reflect the next empty entry point.
$ws.Range("B17").Select() | Out-Null
